$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p007v_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p007v_1</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p007v_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p007v_2</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p007v_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p007v_3</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p007v_a4</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p007v_4</id>", 2) | Out-Null
